$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(26,1).Value2 = 111273661
$ws.Cells.Item(26,2).Value2 = 89686
$ws.Cells.Item(26,5).Value2 = 658
$ws.Cells.Item(26,6).Value2 = 'Rosenticka'
$ws.Cells.Item(26,7).Value2 = 'Rhodofomes roseus'
$ws.Cells.Item(26,8).Value2 = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(26,17).Value2 = 591636.9769660851
$ws.Cells.Item(26,18).Value2 = 7043422.612332962
$ws.Cells.Item(27,1).Value2 = 111273656
$ws.Cells.Item(27,2).Value2 = 73696
$ws.Cells.Item(27,4).Value2 = 'NT'
$ws.Cells.Item(27,5).Value2 = 6440
$ws.Cells.Item(27,6).Value2 = 'Vitgrynig nållav'
$ws.Cells.Item(27,7).Value2 = 'Chaenotheca subroscida'
$ws.Cells.Item(27,8).Value2 = '(Eitner) Zahlbr.'
$ws.Cells.Item(27,17).Value2 = 591725.0424782543
$ws.Cells.Item(27,18).Value2 = 7043424.7006835
$ws.Cells.Item(28,1).Value2 = 111273663
$ws.Cells.Item(28,2).Value2 = 89686
$ws.Cells.Item(28,5).Value2 = 658
$ws.Cells.Item(28,6).Value2 = 'Rosenticka'
$ws.Cells.Item(28,7).Value2 = 'Rhodofomes roseus'
$ws.Cells.Item(28,8).Value2 = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(28,17).Value2 = 591652.4436271309
$ws.Cells.Item(28,18).Value2 = 7043413.675855185
$ws.Cells.Item(29,1).Value2 = 111273670
$ws.Cells.Item(29,2).Value2 = 77515
$ws.Cells.Item(29,4).Value2 = 'NT'
$ws.Cells.Item(29,5).Value2 = 6425
$ws.Cells.Item(29,6).Value2 = 'Garnlav'
$ws.Cells.Item(29,7).Value2 = 'Alectoria sarmentosa'
$ws.Cells.Item(29,8).Value2 = '(Ach.) Ach.'
$ws.Cells.Item(29,17).Value2 = 591622.4606337334
$ws.Cells.Item(29,18).Value2 = 7043398.517451782
$ws.Cells.Item(30,1).Value2 = 111273659
$ws.Cells.Item(30,2).Value2 = 89845
$ws.Cells.Item(30,4).Value2 = 'VU'
$ws.Cells.Item(30,5).Value2 = 1209
$ws.Cells.Item(30,6).Value2 = 'Rynkskinn'
$ws.Cells.Item(30,7).Value2 = 'Phlebia centrifuga'
$ws.Cells.Item(30,8).Value2 = 'P.Karst.'
$ws.Cells.Item(30,17).Value2 = 591495.2093399345
$ws.Cells.Item(30,18).Value2 = 7043327.847347787
$ws.Cells.Item(31,1).Value2 = 111273666
$ws.Cells.Item(31,2).Value2 = 96348
$ws.Cells.Item(31,5).Value2 = 220787
$ws.Cells.Item(31,6).Value2 = 'Knärot'
$ws.Cells.Item(31,7).Value2 = 'Goodyera repens'
$ws.Cells.Item(31,8).Value2 = '(L.) R. Br.'
$ws.Cells.Item(31,17).Value2 = 591499.5271172373
$ws.Cells.Item(31,18).Value2 = 7043317.696102448
$ws.Cells.Item(32,1).Value2 = 111273655
$ws.Cells.Item(32,2).Value2 = 73696
$ws.Cells.Item(32,5).Value2 = 6440
$ws.Cells.Item(32,6).Value2 = 'Vitgrynig nållav'
$ws.Cells.Item(32,7).Value2 = 'Chaenotheca subroscida'
$ws.Cells.Item(32,8).Value2 = '(Eitner) Zahlbr.'
$ws.Cells.Item(32,17).Value2 = 591622.4606337334
$ws.Cells.Item(32,18).Value2 = 7043398.517451782
$ws.Cells.Item(33,1).Value2 = 111273672
$ws.Cells.Item(33,2).Value2 = 77515
$ws.Cells.Item(33,5).Value2 = 6425
$ws.Cells.Item(33,6).Value2 = 'Garnlav'
$ws.Cells.Item(33,7).Value2 = 'Alectoria sarmentosa'
$ws.Cells.Item(33,8).Value2 = '(Ach.) Ach.'
$ws.Cells.Item(33,17).Value2 = 591719.3732997013
$ws.Cells.Item(33,18).Value2 = 7043419.6232786
$ws.Cells.Item(35,1).Value2 = 111273669
$ws.Cells.Item(35,17).Value2 = 591616.806528918
$ws.Cells.Item(35,18).Value2 = 7043377.357856153
$ws.Cells.Item(36,1).Value2 = 111273667
$ws.Cells.Item(36,2).Value2 = 89423
$ws.Cells.Item(36,5).Value2 = 5432
$ws.Cells.Item(36,6).Value2 = 'Granticka'
$ws.Cells.Item(36,7).Value2 = 'Porodaedalea chrysoloma'
$ws.Cells.Item(36,8).Value2 = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(36,17).Value2 = 591618.866522243
$ws.Cells.Item(36,18).Value2 = 7043352.399297187
$ws.Cells.Item(37,1).Value2 = 111273664
$ws.Cells.Item(37,2).Value2 = 89590
$ws.Cells.Item(37,4).Value2 = 'VU'
$ws.Cells.Item(37,5).Value2 = 48
$ws.Cells.Item(37,6).Value2 = 'Lappticka'
$ws.Cells.Item(37,7).Value2 = 'Amylocystis lapponica'
$ws.Cells.Item(37,8).Value2 = '(Romell) Singer'
$ws.Cells.Item(37,17).Value2 = 591673.2841504611
$ws.Cells.Item(37,18).Value2 = 7043420.083276978
$ws.Cells.Item(38,1).Value2 = 111315145
$ws.Cells.Item(38,2).Value2 = 89845
$ws.Cells.Item(38,4).Value2 = 'VU'
$ws.Cells.Item(38,5).Value2 = 1209
$ws.Cells.Item(38,6).Value2 = 'Rynkskinn'
$ws.Cells.Item(38,7).Value2 = 'Phlebia centrifuga'
$ws.Cells.Item(38,8).Value2 = 'P.Karst.'
$ws.Cells.Item(38,17).Value2 = 591478.5830416525
$ws.Cells.Item(38,18).Value2 = 7043314.860723522
$ws.Cells.Item(39,1).Value2 = 111315149
$ws.Cells.Item(39,2).Value2 = 89686
$ws.Cells.Item(39,5).Value2 = 658
$ws.Cells.Item(39,6).Value2 = 'Rosenticka'
$ws.Cells.Item(39,7).Value2 = 'Rhodofomes roseus'
$ws.Cells.Item(39,8).Value2 = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(39,11).ClearContents()
$ws.Cells.Item(39,17).Value2 = 591670.9593730925
$ws.Cells.Item(39,18).Value2 = 7043423.143536596
$ws.Cells.Item(40,1).Value2 = 111315148
$ws.Cells.Item(40,2).Value2 = 89686
$ws.Cells.Item(40,4).Value2 = 'NT'
$ws.Cells.Item(40,5).Value2 = 658
$ws.Cells.Item(40,6).Value2 = 'Rosenticka'
$ws.Cells.Item(40,7).Value2 = 'Rhodofomes roseus'
$ws.Cells.Item(40,8).Value2 = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(40,17).Value2 = 591645.4590963478
$ws.Cells.Item(40,18).Value2 = 7043407.667238996
$ws.Cells.Item(41,1).Value2 = 111268460
$ws.Cells.Item(41,2).Value2 = 89686
$ws.Cells.Item(41,5).Value2 = 658
$ws.Cells.Item(41,6).Value2 = 'Rosenticka'
$ws.Cells.Item(41,7).Value2 = 'Rhodofomes roseus'
$ws.Cells.Item(41,8).Value2 = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(41,17).Value2 = 591472.6953434804
$ws.Cells.Item(41,18).Value2 = 7043317.372138057
$ws.Cells.Item(42,1).Value2 = 111266420
$ws.Cells.Item(42,2).Value2 = 73696
$ws.Cells.Item(42,5).Value2 = 6440
$ws.Cells.Item(42,6).Value2 = 'Vitgrynig nållav'
$ws.Cells.Item(42,7).Value2 = 'Chaenotheca subroscida'
$ws.Cells.Item(42,8).Value2 = '(Eitner) Zahlbr.'
$ws.Cells.Item(42,11).Value2 = ''
$ws.Cells.Item(42,17).Value2 = 591722.3379231346
$ws.Cells.Item(42,18).Value2 = 7043409.880360964
$ws.Cells.Item(43,1).Value2 = 111315146
$ws.Cells.Item(43,2).Value2 = 77515
$ws.Cells.Item(43,4).Value2 = 'NT'
$ws.Cells.Item(43,5).Value2 = 6425
$ws.Cells.Item(43,6).Value2 = 'Garnlav'
$ws.Cells.Item(43,7).Value2 = 'Alectoria sarmentosa'
$ws.Cells.Item(43,8).Value2 = '(Ach.) Ach.'
$ws.Cells.Item(43,9).Value2 = ''
$ws.Cells.Item(43,17).Value2 = 591616.7319226691
$ws.Cells.Item(43,18).Value2 = 7043364.400079632
$ws.Cells.Item(44,1).Value2 = 111267164
$ws.Cells.Item(44,2).Value2 = 89590
$ws.Cells.Item(44,4).Value2 = 'VU'
$ws.Cells.Item(44,5).Value2 = 48
$ws.Cells.Item(44,6).Value2 = 'Lappticka'
$ws.Cells.Item(44,7).Value2 = 'Amylocystis lapponica'
$ws.Cells.Item(44,8).Value2 = '(Romell) Singer'
$ws.Cells.Item(44,11).Value2 = ''
$ws.Cells.Item(44,17).Value2 = 591635.2558426465
$ws.Cells.Item(44,18).Value2 = 7043404.693209249
$ws.Cells.Item(45,1).Value2 = 111268512
$ws.Cells.Item(45,2).Value2 = 56398
$ws.Cells.Item(45,5).Value2 = 100109
$ws.Cells.Item(45,6).Value2 = 'Tretåig hackspett'
$ws.Cells.Item(45,7).Value2 = 'Picoides tridactylus'
$ws.Cells.Item(45,8).Value2 = '(Linnaeus, 1758)'
$ws.Cells.Item(45,11).Value2 = ''
$ws.Cells.Item(45,13).Value2 = 'äldre spår'
$ws.Cells.Item(45,17).Value2 = 591472.6953434804
$ws.Cells.Item(45,18).Value2 = 7043317.372138057
$ws.Cells.Item(46,1).Value2 = 111315143
$ws.Cells.Item(46,2).Value2 = 89686
$ws.Cells.Item(46,4).Value2 = 'NT'
$ws.Cells.Item(46,5).Value2 = 658
$ws.Cells.Item(46,6).Value2 = 'Rosenticka'
$ws.Cells.Item(46,7).Value2 = 'Rhodofomes roseus'
$ws.Cells.Item(46,8).Value2 = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(46,17).Value2 = 591477.5224061215
$ws.Cells.Item(46,18).Value2 = 7043320.638036993
$ws.Cells.Item(47,1).Value2 = 111315139
$ws.Cells.Item(47,2).Value2 = 96348
$ws.Cells.Item(47,4).Value2 = 'VU'
$ws.Cells.Item(47,5).Value2 = 220787
$ws.Cells.Item(47,6).Value2 = 'Knärot'
$ws.Cells.Item(47,7).Value2 = 'Goodyera repens'
$ws.Cells.Item(47,8).Value2 = '(L.) R. Br.'
$ws.Cells.Item(47,9).Value2 = '1'
$ws.Cells.Item(47,17).Value2 = 591510.9235177813
$ws.Cells.Item(47,18).Value2 = 7043279.155835367
$ws.Cells.Item(47,29).Value2 = 'Plus massor av bladrosetter'
$ws.Cells.Item(48,1).Value2 = 111315141
$ws.Cells.Item(48,2).Value2 = 96348
$ws.Cells.Item(48,4).Value2 = 'VU'
$ws.Cells.Item(48,5).Value2 = 220787
$ws.Cells.Item(48,6).Value2 = 'Knärot'
$ws.Cells.Item(48,7).Value2 = 'Goodyera repens'
$ws.Cells.Item(48,8).Value2 = '(L.) R. Br.'
$ws.Cells.Item(48,9).Value2 = '3'
$ws.Cells.Item(48,11).ClearContents()
$ws.Cells.Item(48,17).Value2 = 591486.5005135566
$ws.Cells.Item(48,18).Value2 = 7043319.555657836
$ws.Cells.Item(49,1).Value2 = 111315151
$ws.Cells.Item(49,2).Value2 = 89590
$ws.Cells.Item(49,4).Value2 = 'VU'
$ws.Cells.Item(49,5).Value2 = 48
$ws.Cells.Item(49,6).Value2 = 'Lappticka'
$ws.Cells.Item(49,7).Value2 = 'Amylocystis lapponica'
$ws.Cells.Item(49,8).Value2 = '(Romell) Singer'
$ws.Cells.Item(49,11).ClearContents()
$ws.Cells.Item(49,13).ClearContents()
$ws.Cells.Item(49,17).Value2 = 591670.9593730925
$ws.Cells.Item(49,18).Value2 = 7043423.143536596
$ws.Cells.Item(50,1).Value2 = 111315147
$ws.Cells.Item(50,2).Value2 = 73696
$ws.Cells.Item(50,5).Value2 = 6440
$ws.Cells.Item(50,6).Value2 = 'Vitgrynig nållav'
$ws.Cells.Item(50,7).Value2 = 'Chaenotheca subroscida'
$ws.Cells.Item(50,8).Value2 = '(Eitner) Zahlbr.'
$ws.Cells.Item(50,17).Value2 = 591620.5314988887
$ws.Cells.Item(50,18).Value2 = 7043403.376114395
$ws.Cells.Item(51,1).Value2 = 111266309
$ws.Cells.Item(51,2).Value2 = 77515
$ws.Cells.Item(51,5).Value2 = 6425
$ws.Cells.Item(51,6).Value2 = 'Garnlav'
$ws.Cells.Item(51,7).Value2 = 'Alectoria sarmentosa'
$ws.Cells.Item(51,8).Value2 = '(Ach.) Ach.'
$ws.Cells.Item(51,11).Value2 = ''
$ws.Cells.Item(51,17).Value2 = 591747.0822552936
$ws.Cells.Item(51,18).Value2 = 7043436.057239689
$ws.Cells.Item(52,1).Value2 = 111315150
$ws.Cells.Item(52,2).Value2 = 89369
$ws.Cells.Item(52,4).Value2 = 'LC'
$ws.Cells.Item(52,5).Value2 = 5447
$ws.Cells.Item(52,6).Value2 = 'Vedticka'
$ws.Cells.Item(52,7).Value2 = 'Fuscoporia viticola'
$ws.Cells.Item(52,8).Value2 = '(Schwein.) Murrill'
$ws.Cells.Item(52,11).ClearContents()
$ws.Cells.Item(52,17).Value2 = 591671.190636521
$ws.Cells.Item(52,18).Value2 = 7043415.108879722
$ws.Cells.Item(53,1).Value2 = 111315142
$ws.Cells.Item(53,2).Value2 = 89405
$ws.Cells.Item(53,4).Value2 = 'NT'
$ws.Cells.Item(53,5).Value2 = 1202
$ws.Cells.Item(53,6).Value2 = 'Ullticka'
$ws.Cells.Item(53,7).Value2 = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(53,8).Value2 = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(53,9).Value2 = ''
$ws.Cells.Item(53,17).Value2 = 591469.6177441666
$ws.Cells.Item(53,18).Value2 = 7043315.49674286
$ws.Cells.Item(53,29).ClearContents()
